# Diagrama-de-Gantt.xlsx edit script
# Adds a new activity "Ejecutar tests." to the Gantt schedule (row 31),
# shifting the remaining activities down by one row, updates the
# duration of "Exportar Logic Apps." (now row 32) from 7 to 13 days,
# and refreshes the dependent chart (series ranges already follow the
# worksheet automatically; layout/axis/position are adjusted to match).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Gantt")

# ---------------------------------------------------------------
# 1. Insert a new row at 31 and seed it with the same formatting as
#    the row above it (row 30), then fill in its content.
# ---------------------------------------------------------------
$ws.Rows("31").Insert()

$ws.Range("B30:E30").Copy()
$ws.Range("B31:E31").PasteSpecial(-4122)

$ws.Range("B31").Value = "Ejecutar tests."
$ws.Range("C31").Formula = '=$E30'
$ws.Range("D31").Value = 1
$ws.Range("E31").Formula = '=$C31+$D31'

# ---------------------------------------------------------------
# 2. Fix up the rows that got pushed down one position. The row
#    insert shifts values/styles correctly but this host does not
#    rewrite formula text automatically, so the chained C/E formulas
#    are re-entered explicitly for rows 32-37.
# ---------------------------------------------------------------
$ws.Range("D32").Value = 13

$ws.Range("C32").Formula = '=$E31'
$ws.Range("E32").Formula = '=$C32+$D32'

$ws.Range("C33").Formula = '=$E32'
$ws.Range("E33").Formula = '=$C33+$D33'

$ws.Range("C34").Formula = '=$E33'
$ws.Range("E34").Formula = '=$C34+$D34'

$ws.Range("C35").Formula = '=$E34'
$ws.Range("E35").Formula = '=$C35+$D35'

$ws.Range("C36").Formula = '=$E35'
$ws.Range("E36").Formula = '=$C36+$D36'

$ws.Range("C37").Formula = '=$E36'
$ws.Range("E37").Formula = '=$C37+$D37'

# ---------------------------------------------------------------
# 3. Restore the selected cell shown in the saved view.
# ---------------------------------------------------------------
$ws.Range("X27").Select()

# ---------------------------------------------------------------
# 4. Chart adjustments (the schedule now spans one extra row, so the
#    chart grew taller and was nudged/resized on the sheet; the
#    title/plot-area manual layouts and the value-axis max moved
#    along with it).
# ---------------------------------------------------------------
$chart = $ws.ChartObjects(1).Chart

$title = $chart.ChartTitle
$title.Left = 0.4472062915163405
$title.Top = 0.015860299720599443

$pa = $chart.PlotArea
$pa.InsideLeft = 0.23901110139010401
$pa.InsideTop = 0.1819698352923276
$pa.InsideWidth = 0.73926050354816764
$pa.InsideHeight = 0.75184525847312567

$valueAxis = $chart.Axes(2)
$valueAxis.MaximumScale = 45663

$co = $ws.ChartObjects(1)
$co.Left = 673.1739068651575
$co.Top = 127.5
$co.Width = 725.9999212598425
$co.Height = 465
